# Weekly price update: insert a new most-recent-week row at row 7
# (pushing existing rows 7-44 down to 8-45) and populate it with the
# new week's Arveja Verde price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 7; this shifts rows 7:44 -> 8:45
$ws.Rows.Item(7).Insert()

# Populate the new row 7 with this week's data
$ws.Cells.Item(7, 1).Value2 = 3
$ws.Cells.Item(7, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(7, 3).Value2 = "Coquimbo"
$ws.Cells.Item(7, 4).Value2 = 44537
$ws.Cells.Item(7, 5).Value2 = 5
$ws.Cells.Item(7, 6).Value2 = 100112022
$ws.Cells.Item(7, 7).Value2 = "Arveja Verde"
$ws.Cells.Item(7, 8).Value2 = "Perfection"
$ws.Cells.Item(7, 9).Value2 = "Primera"
$ws.Cells.Item(7, 10).Value2 = 78
$ws.Cells.Item(7, 11).Value2 = 27000
$ws.Cells.Item(7, 12).Value2 = 28000
$ws.Cells.Item(7, 13).Value2 = 27487
$ws.Cells.Item(7, 14).Value2 = "`$/malla 25 kilos"
$ws.Cells.Item(7, 15).Value2 = "Provincia de Limarí"
$ws.Cells.Item(7, 16).Value2 = 1099
$ws.Cells.Item(7, 17).Value2 = 25
$ws.Cells.Item(7, 18).Value2 = "Hortaliza"
